$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("grandes regiões e unidades da federação" header) is removed; the
# rows below (7..37) shift up by one so each region's label now lines up
# with the data that originally belonged to the next region. This is a
# straightforward row deletion of row 6.
$ws.Rows.Item(6).Delete()
